$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 121
$ws.Range("H121").Value = 778.6
$ws.Range("I121").Value = 1097.5
$ws.Range("J121").Value = 755.8214
$ws.Range("K121").Value = 3292.5
$ws.Range("L121").Value = 2267.4642
$ws.Range("M121").Value = -1545.5
$ws.Range("N121").Value = -5761.4642

# Row 134
$ws.Range("H134").Value = 127686.664
$ws.Range("J134").Value = 127686.664
$ws.Range("L134").Value = 127686.664
$ws.Range("N134").Value = -137826.664

# Row 136
$ws.Range("H136").Value = 49834
$ws.Range("J136").Value = 49834
$ws.Range("L136").Value = 49834
$ws.Range("N136").Value = -60034

# Row 137
$ws.Range("H137").Value = 1358.6038
$ws.Range("I137").Value = 1221.4
$ws.Range("K137").Value = 3664.2
$ws.Range("M137").Value = -1114.2

# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("N139").Value = 0

# Row 140
$ws.Range("H140").Value = 81227
$ws.Range("J140").Value = 79857.78
$ws.Range("L140").Value = 79857.78
$ws.Range("N140").Value = -90217.78

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 10731.671
$ws.Range("I32").Value = 11351.698
$ws.Range("J32").Value = 5702.5557
$ws.Range("K32").Value = 11351.698
$ws.Range("L32").Value = 5702.5557
$ws.Range("M32").Value = -11064.698
$ws.Range("N32").Value = -6276.5557

# Row 133
$ws.Range("H133").Value = 41738
$ws.Range("J133").Value = 41738
$ws.Range("L133").Value = 41738
$ws.Range("N133").Value = -46798

# Row 134
$ws.Range("H134").Value = 65167.5
$ws.Range("J134").Value = 65167.5
$ws.Range("L134").Value = 65167.5
$ws.Range("N134").Value = -75307.5

# Row 138
$ws.Range("H138").Value = 61004.46
$ws.Range("J138").Value = 61004.46
$ws.Range("L138").Value = 61004.46
$ws.Range("N138").Value = -71284.45999999999

# Row 139
$ws.Range("H139").Value = 86183.91
$ws.Range("J139").Value = 86183.91
$ws.Range("L139").Value = 86183.91
$ws.Range("N139").Value = -96463.91

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1181.7742
$ws.Range("I94").Value = 1002.5
$ws.Range("J94").Value = 1620
$ws.Range("K94").Value = 1002.5
$ws.Range("L94").Value = 1620
$ws.Range("M94").Value = -551.5
$ws.Range("N94").Value = -2522

# Row 132
$ws.Range("H132").Value = 75347
$ws.Range("J132").Value = 75347
$ws.Range("L132").Value = 75347
$ws.Range("N132").Value = -85467

# Row 133
$ws.Range("H133").Value = 61150
$ws.Range("J133").Value = 61150
$ws.Range("L133").Value = 61150
$ws.Range("N133").Value = -71270

# Row 135
$ws.Range("H135").Value = 57193.332
$ws.Range("J135").Value = 57193.332
$ws.Range("L135").Value = 57193.332
$ws.Range("N135").Value = -67333.33199999999

# Row 138
$ws.Range("H138").Value = 39997.8
$ws.Range("J138").Value = 39997.8
$ws.Range("L138").Value = 39997.8
$ws.Range("N138").Value = -50277.8

# Row 140
$ws.Range("H140").Value = 55512.375
$ws.Range("J140").Value = 55512.375
$ws.Range("L140").Value = 55512.375
$ws.Range("N140").Value = -65872.375

$ws = $wb.Worksheets.Item("CRP")
# Row 133
$ws.Range("H133").Value = 67219
$ws.Range("J133").Value = 67219
$ws.Range("L133").Value = 67219
$ws.Range("N133").Value = -72279

# Row 135
$ws.Range("H135").Value = 94722.78
$ws.Range("J135").Value = 94722.78
$ws.Range("L135").Value = 94722.78
$ws.Range("N135").Value = -104862.78

# Row 137
$ws.Range("H137").Value = 44767.777
$ws.Range("J137").Value = 56582
$ws.Range("L137").Value = 56582
$ws.Range("N137").Value = -66782

# Row 138
$ws.Range("H138").Value = 52046.812
$ws.Range("J138").Value = 52046.812
$ws.Range("L138").Value = 52046.812
$ws.Range("N138").Value = -62326.812

# Row 140
$ws.Range("H140").Value = 80769.914
$ws.Range("J140").Value = 80769.914
$ws.Range("L140").Value = 80769.914
$ws.Range("N140").Value = -91129.914

# Row 141
$ws.Range("H141").Value = 32395.857
$ws.Range("J141").Value = 33628.668
$ws.Range("L141").Value = 33628.668
$ws.Range("N141").Value = -43988.668

$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 520.1724
$ws.Range("I107").Value = 741.1818
$ws.Range("J107").Value = 385.1111
$ws.Range("K107").Value = 2223.5454
$ws.Range("L107").Value = 1155.3333
$ws.Range("M107").Value = -303.5454
$ws.Range("N107").Value = -4995.3333

# Row 113
$ws.Range("H113").Value = 760.7222
$ws.Range("I113").Value = 531.6667
$ws.Range("K113").Value = 1595.0001
$ws.Range("M113").Value = 574.9999

$ws = $wb.Worksheets.Item("GSM")
# Row 123
$ws.Range("H123").Value = 29629.584
$ws.Range("J123").Value = 29629.584
$ws.Range("L123").Value = 29629.584
$ws.Range("N123").Value = -34529.584

# Row 133
$ws.Range("H133").Value = 57875
$ws.Range("J133").Value = 57875
$ws.Range("L133").Value = 57875
$ws.Range("N133").Value = -67995

# Row 135
$ws.Range("H135").Value = 48575.3
$ws.Range("J135").Value = 48575.3
$ws.Range("L135").Value = 48575.3
$ws.Range("N135").Value = -58715.3

# Row 138
$ws.Range("H138").Value = 48571.6
$ws.Range("J138").Value = 48571.6
$ws.Range("L138").Value = 48571.6
$ws.Range("N138").Value = -58851.6

# Row 139
$ws.Range("H139").Value = 115732
$ws.Range("J139").Value = 115732
$ws.Range("L139").Value = 115732
$ws.Range("N139").Value = -126012

# Row 140
$ws.Range("H140").Value = 46846.363
$ws.Range("J140").Value = 46846.363
$ws.Range("L140").Value = 46846.363
$ws.Range("N140").Value = -57206.363

# Row 141
$ws.Range("H141").Value = 44748.332
$ws.Range("J141").Value = 44748.332
$ws.Range("L141").Value = 44748.332
$ws.Range("N141").Value = -55108.332

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3950.2144
$ws.Range("I7").Value = 3472
$ws.Range("J7").Value = 4428.4287
$ws.Range("K7").Value = 3472
$ws.Range("L7").Value = 4428.4287
$ws.Range("M7").Value = -3360
$ws.Range("N7").Value = -4652.4287

# Row 11
$ws.Range("H11").Value = 50000
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

# Row 126
$ws.Range("H126").Value = 3950.2144
$ws.Range("I126").Value = 3472
$ws.Range("J126").Value = 4428.4287
$ws.Range("K126").Value = 10416
$ws.Range("L126").Value = 13285.2861
$ws.Range("M126").Value = -7946
$ws.Range("N126").Value = -18225.2861

# Row 134
$ws.Range("H134").Value = 77425.60000000001
$ws.Range("J134").Value = 77425.60000000001
$ws.Range("L134").Value = 77425.60000000001
$ws.Range("N134").Value = -87565.60000000001

# Row 135
$ws.Range("H135").Value = 58431.332
$ws.Range("J135").Value = 58431.332
$ws.Range("L135").Value = 58431.332
$ws.Range("N135").Value = -68571.33199999999

# Row 137
$ws.Range("H137").Value = 85275
$ws.Range("J137").Value = 85275
$ws.Range("L137").Value = 85275
$ws.Range("N137").Value = -95475

# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").ClearContents()
$ws.Range("N138").Value = 0

# Row 139
$ws.Range("H139").Value = 37879
$ws.Range("J139").Value = 37879
$ws.Range("L139").Value = 37879
$ws.Range("N139").Value = -48159

# Row 140
$ws.Range("H140").Value = 61870.75
$ws.Range("J140").Value = 61870.75
$ws.Range("L140").Value = 61870.75
$ws.Range("N140").Value = -72230.75

# Row 141
$ws.Range("H141").Value = 44156.875
$ws.Range("J141").Value = 44156.875
$ws.Range("L141").Value = 44156.875
$ws.Range("N141").Value = -54516.875

$ws = $wb.Worksheets.Item("WVR")
# Row 8
$ws.Range("H8").Value = 734.3333
$ws.Range("I8").Value = 734.3333
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 734.3333
$ws.Range("L8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -594.3333

# Row 133
$ws.Range("H133").Value = 48215.2
$ws.Range("J133").Value = 48215.2
$ws.Range("L133").Value = 48215.2
$ws.Range("N133").Value = -58335.2

# Row 135
$ws.Range("H135").Value = 54089.918
$ws.Range("J135").Value = 54089.918
$ws.Range("L135").Value = 54089.918
$ws.Range("N135").Value = -64229.918

# Row 136
$ws.Range("H136").Value = 1272
$ws.Range("I136").Value = 1224.75
$ws.Range("K136").Value = 3674.25
$ws.Range("M136").Value = -1124.25

# Row 137
$ws.Range("H137").Value = 58393.8
$ws.Range("J137").Value = 58393.8
$ws.Range("L137").Value = 58393.8
$ws.Range("N137").Value = -68593.8

# Row 139
$ws.Range("H139").Value = 57525
$ws.Range("J139").Value = 57525
$ws.Range("L139").Value = 57525
$ws.Range("N139").Value = -67805

# Row 140
$ws.Range("H140").Value = 34438.875
$ws.Range("J140").Value = 34438.875
$ws.Range("L140").Value = 34438.875
$ws.Range("N140").Value = -44798.875

# Row 141
$ws.Range("H141").Value = 81289.5
$ws.Range("J141").Value = 81289.5
$ws.Range("L141").Value = 81289.5
$ws.Range("N141").Value = -91649.5
